# ----------------------------------------------------------------------------
# "Fruta / hortaliza, semanal" - weekly refresh of the Cebollin/La Palmera sheet.
# A new daily price record is inserted at the front of the date series (row 67),
# which pushes every subsequent record (rows 67-145) down by one row; the record
# that used to be last (row 145) lands in a brand-new row 146.
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append row 146: clone the non-shifting ("static") columns straight from row 145 -
#    every record in this sheet shares the same market/category/unit/origin metadata.
$staticCols = @("A","B","C","E","F","G","H","I","N","O","Q","R")
foreach ($col in $staticCols) {
    $ws.Range($col + "146").Value = $ws.Range($col + "145").Value2
}
$ws.Range("D146").NumberFormat = $ws.Range("D145").NumberFormat

# 2) Shift the per-record columns (Fecha/Volumen/Precio min/max/ponderado/$-Kg) down
#    one row at a time: new row N takes the values the OLD row (N-1) held, for every
#    N from 146 down to 67. Row 67 alone gets a freshly reported date (44557); its
#    other figures come from the old row 66 (which itself is left untouched).
$shiftedValues = @{
    67 = @{ D = 44557; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    68 = @{ D = 44384; J = 3320; K = 900; L = 1000; M = 950; P = 158 }
    69 = @{ D = 44487; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    70 = @{ D = 44326; J = 2680; K = 900; L = 1000; M = 950; P = 158 }
    71 = @{ D = 44407; J = 3200; K = 900; L = 1000; M = 950; P = 158 }
    72 = @{ D = 44354; J = 2600; K = 900; L = 1000; M = 950; P = 158 }
    73 = @{ D = 44505; J = 3100; K = 900; L = 1000; M = 950; P = 158 }
    74 = @{ D = 44274; J = 2700; K = 900; L = 1000; M = 950; P = 158 }
    75 = @{ D = 44309; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    76 = @{ D = 44214; J = 2600; K = 1000; L = 1100; M = 1050; P = 175 }
    77 = @{ D = 44172; J = 2400; K = 800; L = 1000; M = 900; P = 150 }
    78 = @{ D = 44237; J = 2800; K = 1000; L = 1100; M = 1050; P = 175 }
    79 = @{ D = 44312; J = 2700; K = 900; L = 1000; M = 950; P = 158 }
    80 = @{ D = 44162; J = 2400; K = 800; L = 1000; M = 900; P = 150 }
    81 = @{ D = 44349; J = 3300; K = 900; L = 1000; M = 950; P = 158 }
    82 = @{ D = 44554; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    83 = @{ D = 44246; J = 2700; K = 900; L = 1000; M = 950; P = 158 }
    84 = @{ D = 44316; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    85 = @{ D = 44174; J = 2600; K = 800; L = 1000; M = 900; P = 150 }
    86 = @{ D = 44277; J = 2400; K = 900; L = 1000; M = 950; P = 158 }
    87 = @{ D = 44482; J = 3100; K = 900; L = 1000; M = 950; P = 158 }
    88 = @{ D = 44265; J = 3200; K = 900; L = 1000; M = 950; P = 158 }
    89 = @{ D = 44330; J = 2960; K = 900; L = 1000; M = 950; P = 158 }
    90 = @{ D = 44323; J = 2880; K = 900; L = 1000; M = 950; P = 158 }
    91 = @{ D = 44209; J = 2700; K = 1000; L = 1100; M = 1050; P = 175 }
    92 = @{ D = 44165; J = 2200; K = 800; L = 1000; M = 900; P = 150 }
    93 = @{ D = 44267; J = 2400; K = 900; L = 1000; M = 950; P = 158 }
    94 = @{ D = 44263; J = 2600; K = 900; L = 1000; M = 950; P = 158 }
    95 = @{ D = 44533; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    96 = @{ D = 44176; J = 2400; K = 800; L = 1000; M = 900; P = 150 }
    97 = @{ D = 44370; J = 3300; K = 900; L = 1000; M = 950; P = 158 }
    98 = @{ D = 44475; J = 3100; K = 900; L = 1000; M = 950; P = 158 }
    99 = @{ D = 44508; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    100 = @{ D = 44291; J = 2600; K = 900; L = 1000; M = 950; P = 158 }
    101 = @{ D = 44468; J = 3160; K = 900; L = 1000; M = 950; P = 158 }
    102 = @{ D = 44379; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    103 = @{ D = 44498; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    104 = @{ D = 44489; J = 3200; K = 900; L = 1000; M = 950; P = 158 }
    105 = @{ D = 44449; J = 3080; K = 900; L = 1000; M = 950; P = 158 }
    106 = @{ D = 44526; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    107 = @{ D = 44218; J = 2600; K = 1000; L = 1100; M = 1050; P = 175 }
    108 = @{ D = 44447; J = 3400; K = 900; L = 1000; M = 950; P = 158 }
    109 = @{ D = 44167; J = 2700; K = 800; L = 1000; M = 900; P = 150 }
    110 = @{ D = 44328; J = 3240; K = 900; L = 1000; M = 950; P = 158 }
    111 = @{ D = 44160; J = 2700; K = 800; L = 1000; M = 900; P = 150 }
    112 = @{ D = 44251; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    113 = @{ D = 44286; J = 3200; K = 800; L = 1000; M = 900; P = 150 }
    114 = @{ D = 44279; J = 3000; K = 800; L = 1000; M = 900; P = 150 }
    115 = @{ D = 44491; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    116 = @{ D = 44389; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    117 = @{ D = 44391; J = 3360; K = 900; L = 1000; M = 950; P = 158 }
    118 = @{ D = 44510; J = 3200; K = 900; L = 1000; M = 950; P = 158 }
    119 = @{ D = 44232; J = 2600; K = 1000; L = 1100; M = 1050; P = 175 }
    120 = @{ D = 44386; J = 3100; K = 900; L = 1000; M = 950; P = 158 }
    121 = @{ D = 44519; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    122 = @{ D = 44414; J = 3200; K = 900; L = 1000; M = 950; P = 158 }
    123 = @{ D = 44543; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    124 = @{ D = 44321; J = 3200; K = 900; L = 1000; M = 950; P = 158 }
    125 = @{ D = 44281; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    126 = @{ D = 44529; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    127 = @{ D = 44405; J = 3360; K = 900; L = 1000; M = 950; P = 158 }
    128 = @{ D = 44200; J = 2500; K = 1000; L = 1100; M = 1050; P = 175 }
    129 = @{ D = 44459; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    130 = @{ D = 44258; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    131 = @{ D = 44298; J = 2600; K = 900; L = 1000; M = 950; P = 158 }
    132 = @{ D = 44428; J = 3120; K = 900; L = 1000; M = 950; P = 158 }
    133 = @{ D = 44340; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    134 = @{ D = 44452; J = 3100; K = 900; L = 1000; M = 950; P = 158 }
    135 = @{ D = 44270; J = 2200; K = 900; L = 1000; M = 950; P = 158 }
    136 = @{ D = 44195; J = 2400; K = 1000; L = 1100; M = 1050; P = 175 }
    137 = @{ D = 44473; J = 2900; K = 900; L = 1000; M = 950; P = 158 }
    138 = @{ D = 44398; J = 3360; K = 900; L = 1000; M = 950; P = 158 }
    139 = @{ D = 44536; J = 2700; K = 900; L = 1000; M = 950; P = 158 }
    140 = @{ D = 44302; J = 2800; K = 900; L = 1000; M = 950; P = 158 }
    141 = @{ D = 44239; J = 2600; K = 1000; L = 1100; M = 1050; P = 175 }
    142 = @{ D = 44344; J = 2960; K = 900; L = 1000; M = 950; P = 158 }
    143 = @{ D = 44463; J = 3000; K = 900; L = 1000; M = 950; P = 158 }
    144 = @{ D = 44365; J = 2900; K = 900; L = 1000; M = 950; P = 158 }
    145 = @{ D = 44454; J = 3360; K = 900; L = 1000; M = 950; P = 158 }
    146 = @{ D = 44272; J = 3100; K = 800; L = 1000; M = 900; P = 150 }
}

146..67 | ForEach-Object {
    $r = $_
    $vals = $shiftedValues[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
